$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds the (text) game date for each row. It was recorded as
# "5-9-2011-12" (ambiguous / off-by-one because of how the NBA stats
# source displayed dates) and should read "2012-05-09" instead.
$rng = $ws.Range("BF2:BF31")

# Force text formatting first so Excel doesn't reinterpret the literal
# "2012-05-09" as a serial date value when it is assigned below.
$rng.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "2012-05-09"
}

# Restore the default (unstyled) cell style now that the text is in place,
# matching the original workbook's formatting for this column.
$rng.Style = "Normal"
